$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cookieText = @'
Cookie
1. The goal is to provide extra knowledge in a httprequest to server
2. The agreement is that server issued Cookie in httpresponse and browser saved it locally and pack it with new httprequest afterward
3. Two types of cookies in servlets:
    * Non-persistent ck: valid for single session only. It is removed each time when user closes the browser
    * Persistent ck: valid for multiple session . It is not removed each time when user closes the browser. It is removed only if user logout or signout.
4. How to set cookie to httpresponse (At server side):
    $ ...
    $ Cookie loginCookie = new Cookie("user",user);
    $ loginCookie.setMaxAge(30*60);     //set cookie expiry in 30 mins
    $ response.addCookie(loginCookie);
    $ response.sendRedirect("LoginSuccess.jsp");
    $ ...
5. How to read Cookie (At server side):
    $ ...
    $ Cookie[] cookies = request.getCookies();
    $ if(cookies != null){
    $   for(Cookie cookie : cookies){
    $     if(cookie.getName().equals("user")) userName = cookie.getValue();
    $ ...
'@

$sessionText = @'
Session?
1. HTTP protocol and Web Servers are stateless, what it means is that for web server every request is a new request to process and they can’t identify if it’s coming from client that has been sending request previously. So comes the Session concept.
2. Session vs. Cookie:
    * Session at Server side / Cookie at client (browser) side
    * Session is a concept impl by Cookie / Cookie is physical at client (browser) side
3. How to set Session to httpresponse (At server side):
 $ if(userID.equals(user) && password.equals(pwd)){
 $  HttpSession session = request.getSession();
 $  session.setAttribute("user", "Pankaj");
 $  session.setMaxInactiveInterval(30*60); //setting session to expiry in 30 mins
 $  Cookie userName = new Cookie("user", user);
 $  userName.setMaxAge(30*60);
 $  response.addCookie(userName);
 $  response.sendRedirect("LoginSuccess.jsp");
    $ ...
4. How to read Session info (At server side):
    $ if(session.getAttribute("user") == null){
    $   response.sendRedirect("login.html");
    $ }else{
    $   user = (String) session.getAttribute("user");
    $ }
    $ String userName = null;
    $ String sessionID = null;
    $ Cookie[] cookies = request.getCookies();
    $ if(cookies !=null){
    $   for(Cookie cookie : cookies){
    $       if(cookie.getName().equals("user")) userName = cookie.getValue();
    $       if(cookie.getName().equals("JSESSIONID")) sessionID = cookie.getValue();
    $   }
    $ }
5. The session data is stored on server side, usually in text files in a temporary directory. They can not be accessed from outside. The thing connecting a session to a client browser is the session ID, which is usually stored in a cookie. This ID is, and should be, the only thing about your session that is stored on client side.

'@

# Row 38: Web / Cookie / <cookie detail>
$ws.Cells.Item(38, 1).Value = "Web"
$ws.Cells.Item(38, 2).Value = "Cookie"
$ws.Cells.Item(38, 3).Value = $cookieText
$ws.Rows(38).RowHeight = 33

# Row 39: Web / Session / <session detail>
$ws.Cells.Item(39, 1).Value = "Web"
$ws.Cells.Item(39, 2).Value = "Session"
$ws.Cells.Item(39, 3).Value = $sessionText
$ws.Rows(39).RowHeight = 33

# Match the author's final selection (cursor rests on B43 after entry)
$ws.Range("B43").Select() | Out-Null

Write-Output "done"
